$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5498
$ws1.Range("G3").Value = 35
$ws1.Range("F4").Value = 12131
$ws1.Range("F5").Value = 299
$ws1.Range("F7").Value = 181
$ws1.Range("F8").Value = 326
$ws1.Range("F9").Value = 1104

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F3").Value = 5498
$ws4.Range("G4").Value = 35
$ws4.Range("F6").Value = 12131
$ws4.Range("F7").Value = 299
$ws4.Range("F9").Value = 181
$ws4.Range("F12").Value = 326
$ws4.Range("F13").Value = 1104
